# Added Display class, updated images
#
# The "Charger" sheet currently has two boolean-ish columns, G (IsTopBuyed)
# and H (IsNew), that are no longer used. They are removed and replaced by
# a single new column G, "DisplayClass" (header only - no data rows yet).
# Everything that used to live to the right of H shifts one column to the
# left versus its original position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete IsTopBuyed / IsNew columns (G:H) - this shifts every
# later column two places to the left.
$ws.Range("G1:H1").EntireColumn.Delete()

# Make room again for the single replacement column and give it its header.
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").Value = "DisplayClass"

# Reflect the author's last selection in the sheet view.
$ws.Range("G2").Select()
